# The sheet gained a new "Docentes responsáveis:" (responsible instructors)
# block right after the "Objectives:" rows and before "Programa resumido:".
# That block occupies 3 new rows (old row 12 "Programa resumido:" and
# everything below it shifts down by 3, old A1:C20 -> new A1:C23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 12, pushing the old row 12 ("Programa
# resumido:") and everything after it down to row 15 onward.
$ws.Rows.Item(12).Resize(3).Insert()

# New row 12: label only in column A (no B/C content for this row).
$ws.Range("A12").Value = "Docentes responsáveis:"

# New row 13: first instructor, same text repeated in columns B and C
# (column C mirrors the "changes" column used throughout this sheet).
$ws.Range("B13").Value = "5817330 - Larissa de Freitas"
$ws.Range("C13").Value = "5817330 - Larissa de Freitas"

# New row 14: second instructor, again mirrored in B and C.
$ws.Range("B14").Value = "1506103 - Pedro Carlos de Oliveira"
$ws.Range("C14").Value = "1506103 - Pedro Carlos de Oliveira"

# The row-insert cloned formatting into B12/C12 and A13/A14, but those
# cells must stay genuinely empty (no cell entry at all), matching the
# target layout - row 12 has only A, rows 13-14 have only B/C.
$ws.Range("B12").Clear()
$ws.Range("C12").Clear()
$ws.Range("A13").Clear()
$ws.Range("A14").Clear()
